$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "coralcover + year + site",
    "rugosity + year + site",
    "coralcover + site",
    "year + site",
    "spongecover + year + site",
    "year + site + year*site",
    "site",
    "spongecover + site",
    "rugosity + site",
    "rugosity + year",
    "rugosity + site + site*rugosity",
    "coralcover + site + site*coralcover",
    "rugosity + year + year*rugosity",
    "spongecover + site + site*spongecover",
    "coralcover + year",
    "coralcover + year + year*coralcover",
    "rugosity",
    "coralcover",
    "year",
    "spongecover + year",
    "spongecover + year + year*spongecover",
    "spongecover"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

$ws.Columns.Item(1).AutoFit() | Out-Null
